$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (F:G) shifting the existing "Is Significant" column
# (and its data) from F/F2 to H/H2. xlInsertShiftToRight = -4161.
$ws.Range("F1:G2").Insert(-4161)

# New header cells (inherit the header style from the shift/insert).
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"

# New data cells for row 2.
$ws.Range("F2").Value = "[505  84] ; [707  45]"
$ws.Range("G2").Value = "[532.34004474  56.65995526] ; [679.65995526  72.34004474]"
